$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns for the "meta" metrics right before the existing
# "arrecadado_sucesso" block (currently column G), shifting everything
# from G onward five columns to the right (G:K become the new columns,
# old G:V become L:AA).
$ws.Range("G1:K1").EntireColumn.Insert()

# --- Header row (row 1) ---
# (the inserted cells already inherit the bold/centered/bordered header
# style shared by the rest of row 1, so no extra formatting call is needed)
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# --- Data rows ---
# Row 2 (sub / apoia.se)
$ws.Range("G2").Value = 165199.0578149446
$ws.Range("H2").Value = 1205.832538795216
$ws.Range("I2").Value = 2163.288658625353
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 21176.91783511972

# Row 3 (sub / catarse)
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Apply the same currency number format used by the neighbouring
# "arrecadado_*" columns (now shifted to L:P) to the new meta columns.
$ws.Range("G2:K3").NumberFormat = "R$ #,##0.00"

Write-Output "meta columns inserted"
